# Insert a new price-snapshot column (AP) right before the existing
# "nom" / "url_produit" columns. This shifts "nom" -> AQ and
# "url_produit" -> AR, and the new AP column is filled with the same
# snapshot value as column AO (the most recent prior price snapshot),
# mirroring how this tracking sheet grows a new timestamped column
# each time it's re-scraped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift "nom" (old AP) -> AQ and "url_produit" (old AQ) -> AR by
# inserting a new blank column at AP. Excel's InsertMode carries the
# formatting (style) of the column being pushed right, matching the
# bold/centered header style already on AP1/AQ1.
$ws.Range("AP1").EntireColumn.Insert()

# New header for the freshly inserted column: the timestamp of this
# snapshot.
$ws.Range("AP1").Value2 = "2026-01-29 13:51:12"

$lastRow = 206

for ($r = 2; $r -le $lastRow; $r++) {
    $aoCell = $ws.Cells.Item($r, 41)   # column AO
    $apCell = $ws.Cells.Item($r, 42)   # column AP (new)

    $val = $aoCell.Value2
    if ($val -ne "" -and $val -ne $null) {
        # Carry the latest known price into the new snapshot column.
        $apCell.Value2 = $val
    }
    # Rows where AO has no price (out of stock before this snapshot)
    # are left blank in AP too, matching AO's own blank state.
}
